$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new row at row 5, shifting existing rows 5-14 down to 6-15
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "LP solver (linprog or gurobi)"
$ws.Range("B5").Value = "gurobi"
